# "9th Stab - Cosmetic Changes"
# This report tracks weekly analyst price-target actions per brokerage firm.
# A new reporting week (Jun_15 then Jun_17) is inserted as new leading data
# columns, pushing the previous weeks (Jun_13, Jun_10) one column to the
# right. Firms with no new rating action this week are filled with "UN".
# A handful of firms did get a new Jun_15 action, those cells are written
# with the new text and highlighted using the sheet's existing
# "changed this week" fill (copied from an already-highlighted cell so the
# workbook keeps reusing the same style / fill definition).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at B, shifting the existing Jun_13 column
# (B) to D and the existing Jun_10 column (C) to E.
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).Insert()

# New column headers for the newly inserted weeks.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Default every data row's new B (Jun_17) / C (Jun_15) cell to "UN"
# (unchanged / no new action that week) - matches the rest of the sheet.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Firms that actually received a new rating action on 6/15/2018 - overwrite
# their Jun_15 (column C) cell with the real text.
$jun15 = @{
    7  = "6/15/2018,Raises Target,Equal Weight -> Equal Weight,$235.00 -> $250.00"
    10 = "6/15/2018,Reiterates,Overweight,$290.00"
    14 = "6/15/2018,Raises Target,Buy -> Buy,$250.00 -> $275.00"
    17 = "6/15/2018,Raises Target,Neutral -> Neutral,$235.00 -> $260.00"
    20 = "6/15/2018,Raises Target,,$278.00 -> $292.00"
    21 = "6/15/2018,Raises Target,Market Perform -> Market Perform,$200.00 -> $250.00"
}

foreach ($r in $jun15.Keys) {
    $ws.Cells.Item($r, 3).Value = $jun15[$r]
}

# Most of those updated cells are also highlighted with the sheet's existing
# "changed this week" fill (copied from an already-highlighted cell so the
# workbook keeps reusing the same style / fill definition). Row 10 (Piper
# Jaffray) keeps the plain/default look.
$highlightRows = @(7, 14, 17, 20, 21)
$ws.Range("E6").Copy()
foreach ($r in $highlightRows) {
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
}

# Give the two "old" weekly columns (C, D) the same custom width as the
# long standing Jun_10 column (E) so all the historical columns line up.
$ws.Columns.Item(3).ColumnWidth = 7.14
$ws.Columns.Item(4).ColumnWidth = 7.14

$excel.CutCopyMode = 0
